$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in row 16 with the new hours entry
$ws.Range("A16").Value = 44351
$ws.Range("B16").Value = 6
$ws.Range("D16").Value = "Bug fixing; Research and working on uPlot and data visualization"

# Update selection to D17 as in the final file
$ws.Range("D17").Select()

$wb.Save()
